$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The "IPA100" test case occupied row 22 (TCID=IPA100, plus its Jira id /
# description / runmode values). Removing it from the suite means clearing
# that row's contents; Excel drops the now-unused shared strings from the
# table automatically on save.
$ws.Range("A22:E22").ClearContents()

# The row had a tall custom height (270pt) to fit its old long text; once
# the text is gone, restore the row to its natural/default height.
$ws.Rows.Item(22).AutoFit()

# Reflect where the user ended up after deleting the row's content.
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 21
[void]$ws.Range("A21").Select()
